$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws "D2" "29.151.18"
Set-TextValue $ws "E2" "  +0.10%  "
Set-TextValue $ws "D3" "1.834.49"
Set-TextValue $ws "E3" "  +0.07%  "
Set-TextValue $ws "D4" "0.9992"
Set-TextValue $ws "D5" "241.41"
Set-TextValue $ws "E5" "  +0.67%  "
Set-TextValue $ws "D6" "0.6602"
Set-TextValue $ws "E6" "  -0.48%  "
Set-TextValue $ws "D7" "1.000"
Set-TextValue $ws "D8" "0.07421"
Set-TextValue $ws "E8" "  +1.13%  "
Set-TextValue $ws "D9" "0.2932"
Set-TextValue $ws "E9" "  -0.70%  "
Set-TextValue $ws "E10" "  +0.91%  "
Set-TextValue $ws "D11" "0.07754"
Set-TextValue $ws "E11" "  +0.99%  "
Set-TextValue $ws "D12" "1.814.33"
Set-TextValue $ws "E12" "  -1.37%  "
Set-TextValue $ws "D13" "4.998"
Set-TextValue $ws "E13" "  -0.38%  "
Set-TextValue $ws "D14" "0.6670"
Set-TextValue $ws "E14" "  -1.02%  "
Set-TextValue $ws "D15" "83.29"
Set-TextValue $ws "E15" "  -3.46%  "
Set-TextValue $ws "D16" "6.119"
Set-TextValue $ws "E16" "  +0.27%  "
Set-TextValue $ws "D17" "0.000008563"
Set-TextValue $ws "E17" "  +4.15%  "
Set-TextValue $ws "D18" "29.143.89"
Set-TextValue $ws "E18" "  +0.07%  "
Set-TextValue $ws "D19" "2.073.93"
Set-TextValue $ws "E19" "  -0.27%  "
Set-TextValue $ws "D20" "226.71"
Set-TextValue $ws "E20" "  -0.60%  "
Set-TextValue $ws "D21" "12.46"
Set-TextValue $ws "E21" "  -0.21%  "
Set-TextValue $ws "D22" "1.002"
Set-TextValue $ws "E22" "  +0.20%  "
Set-TextValue $ws "D23" "7.075"
Set-TextValue $ws "E23" "  -2.89%  "
Set-TextValue $ws "D24" "0.9999"
Set-TextValue $ws "E24" "  +0.02%  "
Set-TextValue $ws "D25" "160.31"
Set-TextValue $ws "E25" "  -0.08%  "
Set-TextValue $ws "D26" "8.632"
Set-TextValue $ws "E26" "  -0.27%  "
Set-TextValue $ws "E27" "  -1.01%  "
Set-TextValue $ws "D28" "17.99"
Set-TextValue $ws "E28" "  +0.01%  "
Set-TextValue $ws "E29" "  +1.15%  "
Set-TextValue $ws "D30" "4.107"
Set-TextValue $ws "D31" "4.034"
Set-TextValue $ws "E31" "  -1.51%  "
Set-TextValue $ws "D32" "1.186"
Set-TextValue $ws "E32" "  -0.70%  "
Set-TextValue $ws "D34" "1.867"
Set-TextValue $ws "E34" "  +0.29%  "
Set-TextValue $ws "D35" "0.7378"
Set-TextValue $ws "E35" "  -1.10%  "
Set-TextValue $ws "D36" "1.146"
Set-TextValue $ws "E36" "  +1.54%  "
Set-TextValue $ws "D37" "2.660"
Set-TextValue $ws "E37" "  -0.67%  "
Set-TextValue $ws "D38" "1.301.75"
Set-TextValue $ws "E38" "  -1.46%  "
Set-TextValue $ws "E39" "  -0.50%  "
Set-TextValue $ws "D40" "2.737"
Set-TextValue $ws "E40" "  +0.93%  "
Set-TextValue $ws "D41" "0.9194"
Set-TextValue $ws "E41" "  -0.27%  "
Set-TextValue $ws "D42" "6.036"
Set-TextValue $ws "E42" "  +0.18%  "
Set-TextValue $ws "D43" "0.08332"
Set-TextValue $ws "E43" "  +8.33%  "
Set-TextValue $ws "E44" "  +0.11%  "
Set-TextValue $ws "D45" "101.97"
Set-TextValue $ws "E45" "  -1.25%  "
Set-TextValue $ws "D46" "1.956.95"
Set-TextValue $ws "E46" "  -1.20%  "
Set-TextValue $ws "D47" "0.5130"
Set-TextValue $ws "E47" "  -0.80%  "
Set-TextValue $ws "D48" "63.67"
Set-TextValue $ws "E48" "  +0.48%  "
Set-TextValue $ws "B49" "BabyDogeCoin"
Set-TextValue $ws "C49" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws "D49" "0.00000000120"
Set-TextValue $ws "E49" "  -0.64%  "
Set-TextValue $ws "B50" "RenderToken"
Set-TextValue $ws "C50" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D50" "1.752"
Set-TextValue $ws "E50" "  -0.31%  "
Set-TextValue $ws "B51" "Cronos"
Set-TextValue $ws "C51" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws "D51" "0.05842"
Set-TextValue $ws "E51" "  -1.37%  "
